$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cell D1 with trailing spaces
$ws.Range("D1").Value = "Trailing spaces   "

# Add new data cells D2 and D3
$ws.Range("D2").Value = 123
$ws.Range("D3").Value = 456

# Move the selection to E3 (matches the resulting selection in the diff)
$ws.Range("E3").Select()
